$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1) Insert a new row above the current row 2 ("PIS / domestic" POST
#    parameters row) to hold the newly-tracked
#    "domestic.x-swagger-router-controller" entry that sits in front
#    of it in the updated report.
# ------------------------------------------------------------------
$ws.Rows.Item(2).Insert()

# Pick up the PIS section formatting (column styles + borders) from
# the row right below (the shifted-down original row 2) so the new
# row matches the rest of the section instead of inheriting row 1's
# header formatting.
$ws.Range("A3:F3").Copy()
$ws.Range("A2:F2").PasteSpecial(-4122)
$ws.Rows.Item(2).RowHeight = 45

$ws.Range("A2").Value = "PIS"
$ws.Range("B2").Value = "HIGH"
$ws.Range("C2").Value = "/v2_1_1.1/payments/v2_1_1.1/domestic.x-swagger-router-controller"
$ws.Range("D2").Value = "pis`n...`n"
$ws.Range("E2").Value = ""
$ws.Range("F2").Value = "Removed:`n" + [char]0x2022 + " pis"

# ------------------------------------------------------------------
# 2) The "upload Polish API" section for CAF
#    (getConfirmationOfFunds.x-swagger-router-controller) has been
#    removed entirely, per the commit message "deleted upload polish
#    api section, it is now loaded in the program". After the insert
#    above, that row sits at 69 (was 68 before the insert).
# ------------------------------------------------------------------
$ws.Rows.Item(69).Delete()

# ------------------------------------------------------------------
# 3) The AIS "deleteConsent.post.parameters" row's Status flips from
#    MINOR back to UNCATEGORIZED. Row numbering nets back out to its
#    original position (74) once the single insert above and single
#    delete above offset each other.
# ------------------------------------------------------------------
$ws.Range("B74").Value = "UNCATEGORIZED"
